# The edit re-orders the 12 data rows (rows 2-13) of the "Artfynd" sheet
# into a new order while leaving the header row (row 1) untouched. No cell
# values themselves change content-wise - entire rows are moved to new row
# positions. This is implemented as a permutation of whole rows, carried
# out with a single scratch row (row 1000, far outside the used range) so
# that every row-to-row move is a plain copy with no data loss, even
# though it is a full rearrangement rather than a simple insert/shift.
#
# after-row -> before-row (which original row's content ends up at each
# final row position), derived from the target workbook state:
#   2<-12  3<-13  4<-2  5<-3  6<-4  7<-6  8<-11  9<-5  10<-7  11<-8  12<-9  13<-10

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$scratch = 1000   # scratch row well below the sheet's real data (row 1-13)
$lastCol = "AY"   # last used column in this sheet

function RowRange($r) {
    return $ws.Range("A" + $r + ":" + $lastCol + $r)
}

# Move the contents of row $src into row $dst (row-range sized to the
# sheet's used columns). The destination is cleared first: this engine's
# Range.Copy does not blank out destination cells for which the source
# cell is empty, so without the explicit clear, stale values could leak
# through from whatever used to be in the destination row.
function MoveRow($srcRow, $dstRow) {
    (RowRange $dstRow).ClearContents()
    (RowRange $srcRow).Copy((RowRange $dstRow))
}

# Cycle 1 of the row permutation: 2 -> 12 -> 9 -> 5 -> 3 -> 13 -> 10 -> 7 -> 6 -> 4 -> 2
MoveRow 2 $scratch
MoveRow 12 2
MoveRow 9 12
MoveRow 5 9
MoveRow 3 5
MoveRow 13 3
MoveRow 10 13
MoveRow 7 10
MoveRow 6 7
MoveRow 4 6
MoveRow $scratch 4

# Cycle 2 of the row permutation: 8 -> 11 -> 8
MoveRow 8 $scratch
MoveRow 11 8
MoveRow $scratch 11

# Clean up the scratch row so it leaves no trace in the saved workbook.
(RowRange $scratch).ClearContents()
